# Update cryptocurrency price/volume figures scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.158.52"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.677.36"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'214.31"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'22.95"
$ws.Range("E8").Value = "  +7.83%  "
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.914.31"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "1.681.13"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "'4.20"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "'0.561"
$ws.Range("E15").Value = "  +4.86%  "
$ws.Range("D16").Value = "'66.59"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "27.132.97"
$ws.Range("D18").Value = "'235.58"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "0.0₃0743"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  -3.86%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'4.55"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Value = "'147.85"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "1.547.59"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("D36").Value = "'0.609"
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("D37").Value = "'0.949"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "'69.76"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "1.823.64"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'1.65"
$ws.Range("E47").Value = "  +6.82%  "
$ws.Range("D48").Value = "'89.92"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "0.0₆0110"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "'8.26"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("E51").Value = "  +1.01%  "
